# This edit reorders (permutes) the data rows 2-7 on the active sheet.
# Each destination row receives the original contents of columns
# A, B, E, F, G, H, Q, R from the row indicated in the mapping below
# (those are the only columns whose values actually differ row-to-row;
# every other column already holds the same value in every row, so we
# leave them untouched to avoid any incidental type/format changes).
#
#   new row 2 <- old row 3
#   new row 3 <- old row 5
#   new row 4 <- old row 6
#   new row 5 <- old row 7
#   new row 6 <- old row 2
#   new row 7 <- old row 4

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary row-to-row (1-based column index).
$cols = @(1, 2, 5, 6, 7, 8, 17, 18)   # A, B, E, F, G, H, Q, R

# Snapshot the values of those columns for every source row (2-7)
# BEFORE any writes happen, since rows are both sources and destinations.
$snapshots = @{}
for ($r = 2; $r -le 7; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshots[$r] = $rowVals
}

# Destination row -> source row mapping.
$mapping = @{
    2 = 3
    3 = 5
    4 = 6
    5 = 7
    6 = 2
    7 = 4
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshots[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
